$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "DBD" (DB layout detail sheet) - add a new LandSeq field
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("DBD")

# Update PrimaryKey / Index1 definitions to include the new LandSeq column
$ws1.Range("C3").Value = "ClCode1,ClCode2,ClNo,LandSeq"
$ws1.Range("C5").Value = "ClCode1,ClCode2,ClNo,LandSeq"
$ws1.Range("D5").Value = "ClCode1 ASC,ClCode2 ASC,ClNo ASC,LandSeq ASC"

# Insert a new field row right after "ClNo" (row 12), pushing the rest of
# the field list (Reason, OtherReason, CreateDate, ...) down by one row.
$ws1.Rows(12).Insert() | Out-Null

# The freshly inserted row inherits a generic blank style; copy the
# formatting used by the rest of the field table (row 11, "ClNo") so the
# new row matches the table's look (borders/fill/font).
$ws1.Range("A11:G11").Copy() | Out-Null
$ws1.Range("A12:G12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the new "LandSeq" field row
$ws1.Range("A12").Value = 4
$ws1.Range("B12").Value = "LandSeq"
$ws1.Range("C12").Value = "土地序號"
$ws1.Range("D12").Value = "DECIMAL"
$ws1.Range("E12").Value = 3

# Renumber the SEQ column for every field that was shifted down
$ws1.Range("A13").Value = 5
$ws1.Range("A14").Value = 6
$ws1.Range("A15").Value = 7
$ws1.Range("A16").Value = 8
$ws1.Range("A17").Value = 9
$ws1.Range("A18").Value = 10

$ws1.Range("C13").Select()

# ---------------------------------------------------------------------
# Sheet "DBS" (DB layout search/key sheet) - reflect the new key column
# in the "讀取Key條件" / "其他ORDER條件" definitions
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DBS")

$ws2.Range("B2").Value = "ClCode1 = ,AND ClCode2 = ,AND ClNo = ,AND LandSeq ="
$ws2.Range("C2").Value = "ClCode1 ASC,ClCode2 ASC,ClNo ASC,LandSeq ASC"
$ws2.Range("B3").Value = "ClCode1 = ,AND ClCode2 = ,AND ClNo = ,AND LandSeq ="
$ws2.Range("C3").Value = "ClCode1 ASC,ClCode2 ASC,ClNo ASC,LandSeq ASC"

$ws2.Range("B9").Select()
